# Insert a new data row at row 24 (pushing all subsequent rows down by one,
# which also grows the sheet from A1:T140 to A1:T141) and populate it with a
# new "Mango" price record for 2022-11-16 (serial 44881), sourced from Brasil.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("24:24").Insert()

$ws.Cells.Item(24, 1).Value = 11
$ws.Cells.Item(24, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(24, 3).Value = "Bíobío"
$ws.Cells.Item(24, 4).Value = 44881
$ws.Cells.Item(24, 5).Value = 8
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100108
$ws.Cells.Item(24, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(24, 9).Value = 100108002
$ws.Cells.Item(24, 10).Value = "Mango"
$ws.Cells.Item(24, 11).Value = "Sin especificar"
$ws.Cells.Item(24, 12).Value = "Primera"
$ws.Cells.Item(24, 13).Value = 200
$ws.Cells.Item(24, 14).Value = 7500
$ws.Cells.Item(24, 15).Value = 8000
$ws.Cells.Item(24, 16).Value = 7800
$ws.Cells.Item(24, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(24, 18).Value = "Brasil"
$ws.Cells.Item(24, 19).Value = 1950
$ws.Cells.Item(24, 20).Value = 4

# Make sure the new date cell keeps the same date style/format as the rest
# of column D in this workbook.
$ws.Cells.Item(24, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
